# Apply cryptos list update (prices/volumes) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.268.18"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.11%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.550.12"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.42%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D5').Value = "'605.28"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.56%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'144.16"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.04%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'3.549.75"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.44%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.22%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  +2.11%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -0.67%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.97%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.00%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'4.156.32"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.53%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.32%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'30.11"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.96%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.556.22"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.63%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'66.383.36"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.16%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E19').Value = "'  +4.90%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'6.18"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.83%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'14.80"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -1.47%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'430.71"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +1.09%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.610"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.10%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'79.57"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.96%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'3.696.63"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.56%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.15%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.58%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +1.32%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'9.12"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.91%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'7.93"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -1.36%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.13%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'RenzoRestakedETH"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'3.547.63"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.66%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'EthereumClassic"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'25.39"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.10%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  -2.68%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -5.81%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'7.82"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.72%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.67%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'5.58"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.67%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'176.05"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +2.48%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.0846"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -1.60%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.10%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.887"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.77%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.77%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'45.92"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +1.47%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'1.00"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.06%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'2.51"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +3.85%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -1.13%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'25.17"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -3.25%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.82%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'23.22"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +2.46%  "
$ws.Range('E51').Style = 'Normal'
